# Rename the worksheet so the template now targets "Opcion Multiple"
# (multiple choice) questions instead of "Emparejamiento" (matching)
# questions, per the commit:
# "Soporte para multiples opciones en grupo emparejamiento por medio de
#  importacion Excel"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Preguntas de Opcion Multiple"

# Reflect the selection left by the author on the merged title cell
# (A1:E2) when the sheet was last saved.
$ws.Range("A1:E2").Select()

Write-Host "Sheet renamed to '$($ws.Name)'"
